# The document contains three <id>...</id> markers, each split across 3
# separate runs (opening tag / inner value / closing tag), with the
# inner-value run carrying different (inherited) formatting than the
# surrounding tag runs. The edit collapses each trio of runs into a
# single run (taking on the Courier New "tag" formatting of the first
# run in the match) and strips the stray "a" from the id values:
#   p090r_a3 -> p090r_3
#   p090v_a1 -> p090v_1
#   p090v_a2 -> p090v_2

$d = $word.ActiveDocument

$replacements = @(
    @{ Old = "<id>p090r_a3</id>"; New = "<id>p090r_3</id>" },
    @{ Old = "<id>p090v_a1</id>"; New = "<id>p090v_1</id>" },
    @{ Old = "<id>p090v_a2</id>"; New = "<id>p090v_2</id>" }
)

foreach ($r in $replacements) {
    $rng = $d.Content
    $found = $rng.Find.Execute($r.Old, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if ($found) {
        # Assigning .Text on the found range (which spans the original
        # 3 runs) replaces it with a single merged run that takes on the
        # formatting of the first run in the match (the Courier New tag
        # styling), matching the target XML exactly.
        $rng.Text = $r.New
    }
}
